$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" (C) column from 46076 to 46077 for all data rows (2-33)
$ws.Range("C2:C33").Value = 46077

# 2) Row data (A 389-2023 / A 1782-2024) swap between rows 4 and 6, plus other rows 4-33 reshuffle.
#    Update Beteckning (A), Datum (B) and Area (G) per record for every row whose content moved.
$ws.Range("A4").Value = 'A 1782-2024'
$ws.Range("B4").Value = 45307
$ws.Range("G4").Value = 2.7
$ws.Range("A6").Value = 'A 389-2023'
$ws.Range("B6").Value = 44929
$ws.Range("G6").Value = 2.5
$ws.Range("A8").Value = 'A 32610-2024'
$ws.Range("B8").Value = 45513
$ws.Range("G8").Value = 0.5
$ws.Range("A9").Value = 'A 48974-2023'
$ws.Range("B9").Value = 45209
$ws.Range("G9").Value = 4.5
$ws.Range("A10").Value = 'A 4822-2023'
$ws.Range("B10").Value = 44957
$ws.Range("G10").Value = 2.2
$ws.Range("A11").Value = 'A 635-2023'
$ws.Range("B11").Value = 44930
$ws.Range("G11").Value = 0.5
$ws.Range("A12").Value = 'A 48181-2024'
$ws.Range("B12").Value = 45589
$ws.Range("G12").Value = 0.7
$ws.Range("A13").Value = 'A 5817-2025'
$ws.Range("B13").Value = 45694.74113425926
$ws.Range("G13").Value = 1.2
$ws.Range("A14").Value = 'A 24-2023'
$ws.Range("B14").Value = 44928
$ws.Range("G14").Value = 0.5
$ws.Range("A15").Value = 'A 28260-2023'
$ws.Range("B15").Value = 45099
$ws.Range("G15").Value = 5
$ws.Range("A16").Value = 'A 21572-2023'
$ws.Range("B16").Value = 45063
$ws.Range("G16").Value = 1.7
$ws.Range("A17").Value = 'A 18328-2025'
$ws.Range("B17").Value = 45762
$ws.Range("G17").Value = 1.8
$ws.Range("A18").Value = 'A 4256-2025'
$ws.Range("B18").Value = 45685
$ws.Range("G18").Value = 2
$ws.Range("A19").Value = 'A 4481-2024'
$ws.Range("B19").Value = 45327
$ws.Range("G19").Value = 1
$ws.Range("A20").Value = 'A 11517-2024'
$ws.Range("B20").Value = 45372
$ws.Range("G20").Value = 0.7
$ws.Range("A21").Value = 'A 1531-2022'
$ws.Range("B21").Value = 44573
$ws.Range("G21").Value = 1.6
$ws.Range("A22").Value = 'A 4486-2024'
$ws.Range("B22").Value = 45327
$ws.Range("G22").Value = 0.6
$ws.Range("A23").Value = 'A 18332-2025'
$ws.Range("B23").Value = 45762
$ws.Range("G23").Value = 2.5
$ws.Range("A24").Value = 'A 18434-2023'
$ws.Range("B24").Value = 45042
$ws.Range("G24").Value = 0.7
$ws.Range("A25").Value = 'A 53131-2021'
$ws.Range("B25").Value = 44468
$ws.Range("G25").Value = 1.3
$ws.Range("A26").Value = 'A 4487-2024'
$ws.Range("B26").Value = 45327
$ws.Range("G26").Value = 1.9
$ws.Range("A27").Value = 'A 18327-2025'
$ws.Range("B27").Value = 45762
$ws.Range("G27").Value = 0.6
$ws.Range("A28").Value = 'A 4493-2024'
$ws.Range("B28").Value = 45327
$ws.Range("G28").Value = 1.8
$ws.Range("A29").Value = 'A 7731-2026'
$ws.Range("B29").Value = 46062.52008101852
$ws.Range("G29").Value = 5.9
$ws.Range("A30").Value = 'A 10710-2025'
$ws.Range("B30").Value = 45722
$ws.Range("G30").Value = 1.8
$ws.Range("A31").Value = 'A 7727-2026'
$ws.Range("B31").Value = 46062.50420138889
$ws.Range("G31").Value = 1.9
$ws.Range("A32").Value = 'A 34400-2025'
$ws.Range("B32").Value = 45846.61351851852
$ws.Range("G32").Value = 1.3
$ws.Range("A33").Value = 'A 34401-2025'
$ws.Range("B33").Value = 45846.6140162037
$ws.Range("G33").Value = 2.8

# 3) Hyperlink formulas in row 4 and row 6 must follow their records too (S,T,V,W,X,Y).
$ws.Range("S4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 1782-2024 artfynd.xlsx"", ""A 1782-2024"")"
$ws.Range("T4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 1782-2024 karta.png"", ""A 1782-2024"")"
$ws.Range("V4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 1782-2024 FSC-klagomål.docx"", ""A 1782-2024"")"
$ws.Range("W4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 1782-2024 FSC-klagomål mail.docx"", ""A 1782-2024"")"
$ws.Range("X4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 1782-2024 tillsynsbegäran.docx"", ""A 1782-2024"")"
$ws.Range("Y4").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 1782-2024 tillsynsbegäran mail.docx"", ""A 1782-2024"")"
$ws.Range("S6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/artfynd/A 389-2023 artfynd.xlsx"", ""A 389-2023"")"
$ws.Range("T6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/kartor/A 389-2023 karta.png"", ""A 389-2023"")"
$ws.Range("V6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomål/A 389-2023 FSC-klagomål.docx"", ""A 389-2023"")"
$ws.Range("W6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/klagomålsmail/A 389-2023 FSC-klagomål mail.docx"", ""A 389-2023"")"
$ws.Range("X6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsyn/A 389-2023 tillsynsbegäran.docx"", ""A 389-2023"")"
$ws.Range("Y6").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1278/tillsynsmail/A 389-2023 tillsynsbegäran mail.docx"", ""A 389-2023"")"

# 4) "Kommuner" label (Markägare, column F) moves from row 12 to row 30.
$ws.Range("F12").ClearContents()
$ws.Range("F30").Value = "Kommuner"
